$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") on rows 2-6 was bumped by one day (45174 -> 45175)
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45175
}
